$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (column C) date value from 2024-10-02 (45567)
# to 2024-10-03 (45568) for all data rows (rows 2 through 29).
for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 3).Value = 45568
}
